$wb = $excel.ActiveWorkbook

# Rename sheets: "Air Sample" -> "AirSample", "Soil Sample" -> "SoilSample"
$wb.Worksheets.Item("Air Sample").Name = "AirSample"
$wb.Worksheets.Item("Soil Sample").Name = "SoilSample"

# SampleCollection sheet: remove "id" column (A1), keep only "samples" in A1
$wsSampleCollection = $wb.Worksheets.Item("SampleCollection")
$wsSampleCollection.Range("A1").Value = "samples"
$wsSampleCollection.Range("B1").Value = $null

# Sample sheet: rename E1 "sample biome" -> "sample_biome", add F1 "sample_type"
$wsSample = $wb.Worksheets.Item("Sample")
$wsSample.Range("E1").Value = "sample_biome"
$wsSample.Range("F1").Value = "sample_type"

# AirSample sheet: rename F1 "sample biome" -> "sample_biome", add G1 "sample_type"
$wsAirSample = $wb.Worksheets.Item("AirSample")
$wsAirSample.Range("F1").Value = "sample_biome"
$wsAirSample.Range("G1").Value = "sample_type"

# SoilSample sheet: rename F1 "sample biome" -> "sample_biome", add G1 "sample_type"
$wsSoilSample = $wb.Worksheets.Item("SoilSample")
$wsSoilSample.Range("F1").Value = "sample_biome"
$wsSoilSample.Range("G1").Value = "sample_type"
